$d = $word.ActiveDocument

# Find the paragraph containing the "R1: Register players" requirement and
# remove it entirely (it was finalized / consolidated elsewhere).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "R1:\s*Regist") {
        $p.Range.Delete()
        break
    }
}
